$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.886.37'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.91%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''1.640.20'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -1.57%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.26%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = '''213.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +2.34%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = '''  -0.17%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = '''  -0.23%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = '''0.2597'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -0.13%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''0.06319'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.18%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = '''20.57'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -2.12%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = '''0.07673'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +1.71%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = '''1.642.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -1.66%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = '''4.406'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -0.16%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = '''1.862.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -1.60%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''0.5471'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  +1.53%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = ('''0.0{0}8179' -f [char]0x2085)
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +2.31%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '''64.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -2.66%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''25.884.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -1.06%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = '''  -0.15%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = '''4.677'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -0.77%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''188.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +0.51%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''10.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -1.17%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = '''6.238'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +0.39%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = '''1.002'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -0.29%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''143.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -3.80%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = '''0.1235'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +0.53%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = '''7.349'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -0.97%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = '''15.84'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +0.89%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = '''1.399'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +2.20%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = '''0.05886'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -4.92%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = '''1.255'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -1.28%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = '''  -0.44%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = '''3.386'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -2.72%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = '''1.634'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -0.14%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = '''0.9836'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -0.88%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = '''  +0.03%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = '''2.735'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -0.90%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = '''0.5587'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -5.33%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = '''0.01593'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -0.25%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = '''5.818'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -3.45%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = '''0.8487'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -0.61%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = '''1.002'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -0.20%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = '''1.018.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  -8.06%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = '''98.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -1.53%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = '''1.788.04'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -1.69%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = ('''0.0{0}108' -f [char]0x2088)
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -3.05%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = '''55.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.04%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = '''1.003'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +0.05%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = '''8.008'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -0.41%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = '''0.05141'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -2.18%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = '''0.4209'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -1.19%  '
$ws.Range("E51").Style = "Normal"
